$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in C2:F16 is stored as text (numbers-as-text, same as source file).
# Re-apply Text format to the range so updated values keep that same semantic.
$ws.Range("C2:F16").NumberFormat = "@"

$ws.Range("C2").Value = "43"
$ws.Range("D2").Value = "39"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "1"

$ws.Range("C3").Value = "48"
$ws.Range("D3").Value = "39"
$ws.Range("E3").Value = "3"
$ws.Range("F3").Value = "0"

$ws.Range("C4").Value = "72"
$ws.Range("D4").Value = "53"
$ws.Range("E4").Value = "7"
$ws.Range("F4").Value = "2"

$ws.Range("C5").Value = "50"
$ws.Range("D5").Value = "43"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "1"

$ws.Range("C6").Value = "43"
$ws.Range("D6").Value = "32"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "2"

$ws.Range("C7").Value = "3"
$ws.Range("D7").Value = "11"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "0"

$ws.Range("C8").Value = "29"
$ws.Range("D8").Value = "24"
$ws.Range("E8").Value = "2"
$ws.Range("F8").Value = "1"

$ws.Range("C9").Value = "18"
$ws.Range("D9").Value = "17"
$ws.Range("E9").Value = "2"
$ws.Range("F9").Value = "0"

$ws.Range("C10").Value = "6"
$ws.Range("D10").Value = "7"
$ws.Range("E10").Value = "0"
$ws.Range("F10").Value = "0"

$ws.Range("C11").Value = "7"
$ws.Range("D11").Value = "7"
$ws.Range("E11").Value = "0"
$ws.Range("F11").Value = "0"

$ws.Range("C12").Value = "33"
$ws.Range("D12").Value = "28"
$ws.Range("E12").Value = "1"
$ws.Range("F12").Value = "0"

$ws.Range("C13").Value = "90"
$ws.Range("D13").Value = "52"
$ws.Range("E13").Value = "4"
$ws.Range("F13").Value = "4"

$ws.Range("C14").Value = "14"
$ws.Range("D14").Value = "13"
$ws.Range("E14").Value = "0"
$ws.Range("F14").Value = "0"

$ws.Range("C15").Value = "1"
$ws.Range("D15").Value = "5"
$ws.Range("E15").Value = "0"
$ws.Range("F15").Value = "0"

$ws.Range("C16").Value = "9"
$ws.Range("D16").Value = "14"
$ws.Range("E16").Value = "0"
$ws.Range("F16").Value = "0"
